$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 45; existing rows 45-117 shift down to 46-118
$ws.Rows.Item(45).Insert()

# Populate the new row 45 with the new weekly data point
$ws.Cells.Item(45, 1).Value = 10
$ws.Cells.Item(45, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(45, 3).Value = "La Araucanía"
$ws.Cells.Item(45, 4).Value = 44579
$ws.Cells.Item(45, 5).Value = 9
$ws.Cells.Item(45, 6).Value = 100112012
$ws.Cells.Item(45, 7).Value = "Espinaca"
$ws.Cells.Item(45, 8).Value = "Sin especificar"
$ws.Cells.Item(45, 9).Value = "Primera"
$ws.Cells.Item(45, 10).Value = 40
$ws.Cells.Item(45, 11).Value = 12000
$ws.Cells.Item(45, 12).Value = 12000
$ws.Cells.Item(45, 13).Value = 12000
$ws.Cells.Item(45, 14).Value = "$/docena de atados"
$ws.Cells.Item(45, 15).Value = "Región de La Araucanía"
$ws.Cells.Item(45, 16).Value = 4000
$ws.Cells.Item(45, 17).Value = 3
$ws.Cells.Item(45, 18).Value = "Hortaliza"

# Apply the same date number format (s="2") used by the rest of column D
$ws.Cells.Item(45, 4).NumberFormat = $ws.Cells.Item(46, 4).NumberFormat
